$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data for columns GW:HH (12 new daily columns), rows 2-7,
# extending the existing GV-ending table.
$data = @(
    @(12,7,0,2,0,10,3,10,12,4,0,0),
    @(20,20,12,5,0,17,6,15,13,10,18,4),
    @(23,24,15,6,2,24,8,17,15,19,27,7),
    @(29,29,21,8,5,25,16,22,23,25,28,16),
    @(31,31,28,9,7,27,27,27,28,33,32,21),
    @(35,36,34,24,27,30,30,36,20,34,34,25)
)

# GW is column 205, HH is column 216
$startCol = 205
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $startCol + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}

# Update the view to match the new extent of the table
$ws.Range("GA1").Select() | Out-Null
$ws.Range("HG2:HH7").Select() | Out-Null
